# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.286832544864788, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 0, 4.840633575959121)
    3  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    4  = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 0, 6.348428708163715)
    5  = @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 3.56341032713086)
    6  = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 3.754798637575387)
    7  = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 0, 8.974608811992548)
    8  = @([double]"6.428364285326893e-09", 0.002571899574220771, 0.1494219747398047, 0.4942365360607697, 1, 0.6462304168031594)
    9  = @(0.6606524410359556, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 0, 1.611132179096228)
    10 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 0, 8.974608811992548)
    11 = @(0.1190320826869504, 10.34677158129881, 3.537761648806719, 10.19245300693656, 0, 24.19601831972904)
    12 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
